# "merged jacobs comments into the structure"
#
# For the assessment of Jacob (row 3 of the first "Peer and self assessment"
# block, and row 22 of the second block) fill in the grade (column B) and
# the written comments (column C) that were merged in from Jacob's review.

# Workbook was re-saved with multi-threaded calculation turned off
# (OOXML calcPr concurrentCalc="0").
$excel.MultiThreadedCalculation.Enabled = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")
$ws.Activate()

# --- Row 22 (second assessment block) : Jacob -----------------------------
# Grade = "Sufficient"
$ws.Range("B22").Value = "Sufficient"
# Comments (written by Jacob about Discord activity / communication)
$ws.Range("C22").WrapText = $true
$ws.Range("C22").Value = "1) Active on Discord. `n2) Good english communication and contribution to working environment`nand development of product. "

# --- Row 3 (first assessment block) : Jacob --------------------------------
# Grade = "Good"
$ws.Range("B3").Value = "Good"
# Comments (longer reflection on working environment / collaboration)
$ws.Range("C3").WrapText = $true
$ws.Range("C3").Value = "1) Really good at creating a great working environment where the social`nrelations between members of the project are prioritized and benefit greatly.`n2) Motivated about the project and great at coming up with ideas and `nparticipating in discussions. `n3) Has not been very active in the process of conducting documented work and sharing files and knowledge with other group members via the online platforms."

# --- Update the view so it is scrolled/selected on the newly filled cell --
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("C3").Select()
